$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row 3 data (SNO, LEVEL, Portal, Topic, QNO, QuestionName, Url, Status)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Medium"
$ws.Range("C3").Value = "LeetCode"
$ws.Range("D3").Value = "Array&Maths"
$ws.Range("E3").Value = 1131
$ws.Range("F3").Value = "Maximum of Absolute Value Expression"
$ws.Range("G3").Value = "https://leetcode.com/problems/maximum-of-absolute-value-expression/"
$ws.Range("H3").Value = "Done"

# Widen columns F and G to fit the new, longer content
$ws.Columns.Item(6).ColumnWidth = 38
$ws.Columns.Item(7).ColumnWidth = 64.3

# Move the active selection to I6
$ws.Range("I6").Select() | Out-Null
